# Fruta / hortaliza, semanal
# Insert a new weekly price record as the newest row for the
# "Feria Lagunitas de Puerto Montt - Alcachofa" subset. The new record is
# inserted as row 8 (pushing the existing rows 8-20 down to rows 9-21), so
# that the data stays ordered the way the source file keeps it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 8:20 down to 9:21, leaving a fresh blank row 8 that
# inherits the formatting (incl. the date-cell style) of the row above it.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row with this week's entry.
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C8").Value = "Los Lagos"
$ws.Range("D8").Value = 44775
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 100112013
$ws.Range("G8").Value = "Alcachofa"
$ws.Range("H8").Value = "Madrigal"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 16000
$ws.Range("L8").Value = 17000
$ws.Range("M8").Value = 16500
$ws.Range("N8").Value = "$/caja 40 unidades"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 412
$ws.Range("Q8").Value = 40
$ws.Range("R8").Value = "Hortaliza"
